$d = $word.ActiveDocument

# --- Step 1: merge the three runs describing the "midcap nifty" definition into a single run ---
# (the "Midcap Nifty:" label run and the following space run are left untouched, matching the diff)
$midcapRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.StartsWith("Midcap Nifty:")) {
        $midcapRange = $candidate.Range
        break
    }
}
$midcapRange.InsertXML('<w:p><w:r><w:t>Midcap Nifty:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>An index or collection of stocks that represent mid-sized businesses on the stock market is commonly referred to as "midcap nifty." The top 50 mid-cap companies listed on the National Stock Exchange (NSE) are tracked by the widely used Nifty Midcap 50 index in India.</w:t></w:r></w:p>')

# --- Step 2: replace the third trailing empty paragraph with the new "Page 2" section content ---
$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)
$targetRange = $target.Range
$targetRange.InsertXML('<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>Page 2</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>What is investing?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The act of distributing funds or resources with the hope of making a profit at some point in the future is known as investing. It entails investing in assets with the intention of generating returns through dividends, interest, or capital growth, such as stocks, bonds, property, or mutual funds. The goal of investing is to accumulate wealth over time and meet financial goals such as funding education, saving for retirement, or accumulating wealth. Thorough investigation, evaluation, and risk control are necessary for successful investing in order to reach well-informed choices and long-term financial objectives.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>What is trading?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Buying and selling financial instruments, such as stocks, currencies, commodities, or derivatives, with the intention of profiting from transient price fluctuations is referred to as trading. Trading focuses on taking advantage of short-term market fluctuations, in contrast to investing, which usually entails holding assets for the long term in order to achieve capital appreciation or income generation.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The goal of traders'' transactions is to profit from price changes that occur over short periods of time, from seconds to days. To find opportunities and make wise choices, they might use a variety of trading strategies, such as technical analysis, fundamental analysis, algorithmic trading, and quantitative analysis.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">A thorough understanding of trading strategies, risk management tactics, and market dynamics are necessary for successful trading. It also entails keeping up with business news, geopolitical developments, economic indicators, and other elements that could affect asset values. Trading platforms, like online brokerages or trading software, are frequently used by traders to execute trades and keep an eye on the state of the market in real time. </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:br/></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>')
